$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'2.30%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.81%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.185"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.38%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07470"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.14%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.440"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'46.26%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.010"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'3.13%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9189"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.68%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1741"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.73%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07613"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.48%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08195"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.54%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03009"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.72%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09940"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.51%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001514"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'1.03%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.006080"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-7.15%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'1.86%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.861"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.61%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.227"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.10%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3261"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.1338"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.15%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.657"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.53%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04610"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-1.21%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.1564"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.83%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001262"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'3.25%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004527"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.58%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-7.27%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0002739"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'51.29%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01785"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'0.04563"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.00%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007386"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'5.59%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'1.47%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002179"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'5.68%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-18.66%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006449"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'5.99%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'15.26%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-19.31%"
$ws.Range("E47").Style = "Normal"
